$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New observation rows appended to the "Artfynd" sheet (rows 15-18),
# extending the used range from A1:AY14 to A1:AY18.

$rows = @(
    @{
        row = 15
        A = 111837758; B = 90187; C = "Ovaliderad"; D = "NT"; E = 2014
        F = "Koralltaggsvamp"; G = "Hericium coralloides"; H = "(Scop.:Fr.) Pers."
        I = "6"; J = "fruktkroppar"
        K = ""; N = ""
        P = "Brotorp, hyggeskant, Sm"
        Q = 575673.5681218; R = 6404513.458820416; S = 10
        T = "Kalmar"; U = "Västervik"; V = "Småland"; W = "Hallingeberg"
        Y = "2023-09-01"; Z = "00:00"; AA = "2023-09-01"; AB = "00:00"
        AC = "På asplåga."
        AD = $false; AE = $false; AF = ""; AG = $false
        AT = ""
        AW = "Magnus Kasselstrand"; AX = "Magnus Kasselstrand"
        AY = ""
    },
    @{
        row = 16
        A = 111837705; B = 90662; C = "Ovaliderad"; D = "LC"; E = 4363
        F = "Zontaggsvamp"; G = "Hydnellum concrescens"; H = "(Pers.) Banker"
        I = "10"; J = "fruktkroppar"
        K = ""; N = ""
        P = "Brotorp, Långsjön, Sm"
        Q = 575795.3141537429; R = 6404518.948622406; S = 10
        T = "Kalmar"; U = "Västervik"; V = "Småland"; W = "Hallingeberg"
        Y = "2023-09-01"; Z = "00:00"; AA = "2023-09-01"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AT = ""
        AW = "Magnus Kasselstrand"; AX = "Magnus Kasselstrand"
        AY = ""
    },
    @{
        row = 17
        A = 111837675; B = 103288; C = "Ovaliderad"; D = "LC"; E = 221144
        F = "Grönpyrola"; G = "Pyrola chlorantha"; H = "Sw."
        I = "10"; J = "plantor/tuvor"
        K = ""; L = ""; N = ""
        P = "Brotorp, Långsjön, Sm"
        Q = 575781.9606960951; R = 6404546.96767282; S = 10
        T = "Kalmar"; U = "Västervik"; V = "Småland"; W = "Hallingeberg"
        Y = "2023-09-01"; Z = "00:00"; AA = "2023-09-01"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AT = ""
        AW = "Magnus Kasselstrand"; AX = "Magnus Kasselstrand"
        AY = ""
    },
    @{
        row = 18
        A = 111837741; B = 90658; C = "Ovaliderad"; D = "NT"; E = 4361
        F = "Orange taggsvamp"; G = "Hydnellum aurantiacum"; H = "(Batsch:Fr.) P.Karst."
        I = "15"; J = "fruktkroppar"
        K = ""; N = ""
        P = "Brotorp, hyggeskant, Sm"
        Q = 575653.9215098171; R = 6404506.688862759; S = 10
        T = "Kalmar"; U = "Västervik"; V = "Småland"; W = "Hallingeberg"
        Y = "2023-09-01"; Z = "00:00"; AA = "2023-09-01"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AT = ""
        AW = "Magnus Kasselstrand"; AX = "Magnus Kasselstrand"
        AY = ""
    }
)

# Numeric columns (plain numbers).
$numericCols = @("A","B","E","Q","R","S")
# Boolean columns.
$boolCols = @("AD","AE","AG")
# Columns that must be stored as *text*, even when the content looks like a
# number or a date/time, matching the source data (quote-prefixed entry).
$textForcedCols = @("I","K","L","N","Y","Z","AA","AB","AF","AT","AY")

foreach ($r in $rows) {
    $rowNum = $r.row
    foreach ($col in $r.Keys) {
        if ($col -eq "row") { continue }
        $addr = "$col$rowNum"
        $val = $r[$col]
        if ($numericCols -contains $col) {
            $ws.Range($addr).Value = $val
        } elseif ($boolCols -contains $col) {
            $ws.Range($addr).Value = $val
        } elseif ($textForcedCols -contains $col) {
            $ws.Range($addr).Value = "'" + $val
            $ws.Range($addr).Style = "Normal"
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
